$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") gets a new table
#    style applied (tableStyleId changes from {2BD116F1-...} to {CBCAC56F-...}).
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{CBCAC56F-56A4-4895-AFBA-DED444DACE7B}")

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's theme: swap the "Integral" (Red Violet) palette
#    that is currently applied to the slide master for the classic
#    "Office Theme" palette.
# ---------------------------------------------------------------------------
function HexToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgbLong $officeThemeColors[$i - 1]
}
